# add test self define coordinate type.lua
# Adds a "coordinate" self-defined type test block (columns G:I) to the
# "monster" sheet, widens/narrows a few columns, and flips the active
# sheet/selection from "role" back to "monster".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("monster")
$role = $wb.Worksheets.Item("role")

# ---- new header / sample data in columns G:I ----------------------------
# (cells are written in this exact order so the shared-string table is
# built up identically to the authored workbook)
$ws.Range("G1").Value = "c1"
$ws.Range("H1").Value = "c2"
$ws.Range("G2").Value = "coordinate"
$ws.Range("I1").Value = "c3"
$ws.Range("H2").Value = "array<coordinate>"
$ws.Range("I2").Value = "group<coordinate>"

$ws.Range("H1").VerticalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4108

$ws.Range("G4").Value = "[]"
$ws.Range("G4").VerticalAlignment = -4108

$ws.Range("G5").Value = "[1,2,3]"
$ws.Range("H5").Value = "[1,2,3],[1,2,3,4],[1,2,3]"
$ws.Range("I5").Value = "{[1,2,3],[1,2]},{[1],[2]}"

$ws.Range("I6").Value = "{[1,2,3],[1,2]},{[1],[3]}"
$ws.Range("I7").Value = "{[1,2,3],[1,2]},{[1],[4]}"
$ws.Range("I8").Value = "{[1,2,3],[1,2]},{[1],[5]}"
$ws.Range("I9").Value = "{[1,2,3],[1,2]},{[1],[6]}"
$ws.Range("I10").Value = "{[1,2,3],[1,2]},{[1],[7]}"
$ws.Range("I11").Value = "{[1,2,3],[1,2]},{[1],[8]}"
$ws.Range("I12").Value = "{[1,2,3,4],[1,2]},{[1],[9]}"

# ---- column widths --------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 9.160714285714286
$ws.Columns.Item(7).ColumnWidth = 14.035714285714286
$ws.Columns.Item(8).ColumnWidth = 37.785714285714285
$ws.Columns.Item(9).ColumnWidth = 125.66071428571429

# ---- selection / active sheet ---------------------------------------------
# select the new cell on "role" first so that activating "monster" below
# becomes the final (and only) active tab
$role.Range("B28").Select()

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$ws.Range("I12").Select()
